# Update results sheet: new distance sweep rows, a new "remove month" error
# column, and a random-forest hyper-parameter sweep table (max depth / n
# estimators) further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add a 5th column + widen the header row ---------
$ws.Range("E1").Value = 'error if we remove "month"'
$ws.Range("D1:G1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

# --- Existing distance-sweep rows 2-4: fill in the classifier/error data -
$ws.Range("A2").Value = 16
$ws.Range("B2").Value = "random forest"
$ws.Range("C2").Value = "xx"
$ws.Range("D2").Value = 72

$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "random forest"
$ws.Range("C3").Value = "xx"
$ws.Range("D3").Value = 69.98

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "random forest"
$ws.Range("C4").Value = "xx"
$ws.Range("D4").Value = 68.92

# --- Row 5 gains the new "remove month" figure ---------------------------
$ws.Range("E5").Value = 64.42

# --- New table starting at row 12: random-forest parameter sweep ---------
$ws.Range("A12").Value = "max depth "
$ws.Range("B12").Value = "n estimators"
$ws.Range("C12").Value = "error"
$ws.Range("D12").Value = "REMOVE WEEKDAY AND MONTH"

$ws.Range("A13").Value = 2
$ws.Range("B13").Value = 50
$ws.Range("C13").Value = 75.83
$ws.Range("D13").Value = 58.67
$ws.Range("D13").Font.Bold = $true

$ws.Range("A14").Value = 5
$ws.Range("B14").Value = 50
$ws.Range("C14").Value = 71.67
$ws.Range("C14").NumberFormat = "General"

$ws.Range("A15").Value = 10
$ws.Range("B15").Value = 50
$ws.Range("C15").Value = 66.57

$ws.Range("A16").Value = 12
$ws.Range("B16").Value = 50
$ws.Range("C16").Value = 65.21

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 50
$ws.Range("C17").Value = 64.83
$ws.Range("A17").Font.Bold = $true
$ws.Range("C17").Font.Bold = $true

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 50
$ws.Range("C18").Value = 65.12

$ws.Range("A19").Value = 20
$ws.Range("B19").Value = 50
$ws.Range("C19").Value = 65.86
$ws.Range("A19").Font.Bold = $false

$ws.Range("A20").Value = 30
$ws.Range("B20").Value = 50
$ws.Range("C20").Value = 67.07

$ws.Range("A21").Value = 40
$ws.Range("B21").Value = 50
$ws.Range("C21").Value = 67.18

$ws.Range("A22").Value = 15
$ws.Range("B22").Value = 10
$ws.Range("C22").Value = 66.29

$ws.Range("A23").Value = 15
$ws.Range("B23").Value = 20
$ws.Range("C23").Value = 65.29

$ws.Range("A24").Value = 15
$ws.Range("B24").Value = 40
$ws.Range("C24").Value = 64.77

$ws.Range("A25").Value = 15
$ws.Range("B25").Value = 80
$ws.Range("C25").Value = 64.39

$ws.Range("A26").Value = 15
$ws.Range("B26").Value = 160
$ws.Range("C26").Value = 64.48

# --- Column width for the new G column + selection cursor ---------------
$ws.Columns.Item(7).ColumnWidth = 10.5
$ws.Range("D17").Select()

Write-Host "applied"
